$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '70.451.16'
$ws.Range('E2').Value = '  -0.96%  '
Set-TextValue $ws.Range('D3') '3.767.05'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '694.50'
$ws.Range('E5').Value = '  -1.35%  '
Set-TextValue $ws.Range('D6') '167.48'
$ws.Range('E6').Value = '  -2.57%  '
Set-TextValue $ws.Range('D7') '3.766.42'
$ws.Range('E7').Value = '  -1.70%  '
$ws.Range('E8').Value = '  +0.43%  '
Set-TextValue $ws.Range('D9') '0.518'
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('E11').Value = '  +1.69%  '
Set-TextValue $ws.Range('D12') '0.473'
$ws.Range('E12').Value = '  +3.04%  '
Set-TextValue $ws.Range('D13') '0.0000245'
$ws.Range('E13').Value = '  -3.64%  '
Set-TextValue $ws.Range('D14') '35.60'
$ws.Range('E14').Value = '  -3.05%  '
Set-TextValue $ws.Range('D15') '4.404.79'
$ws.Range('E15').Value = '  -1.63%  '
Set-TextValue $ws.Range('D16') '3.854.65'
$ws.Range('E16').Value = '  +0.36%  '
Set-TextValue $ws.Range('D17') '70.615.63'
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D19') '7.09'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D20') '17.32'
$ws.Range('E20').Value = '  -0.45%  '
Set-TextValue $ws.Range('D21') '513.45'
$ws.Range('E21').Value = '  +3.90%  '
Set-TextValue $ws.Range('D22') '10.28'
$ws.Range('E22').Value = '  -4.04%  '
Set-TextValue $ws.Range('D23') '0.707'
$ws.Range('E23').Value = '  -4.04%  '
Set-TextValue $ws.Range('D24') '83.07'
$ws.Range('E24').Value = '  -2.63%  '
$ws.Range('E25').Value = '  -4.80%  '
Set-TextValue $ws.Range('D26') '12.43'
$ws.Range('E26').Value = '  +2.70%  '
Set-TextValue $ws.Range('D27') '3.917.80'
$ws.Range('E27').Value = '  -1.64%  '
Set-TextValue $ws.Range('D28') '10.06'
$ws.Range('E28').Value = '  -5.38%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D30') '1.92'
$ws.Range('E30').Value = '  -8.10%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D31') '2.91'
$ws.Range('E31').Value = '  -5.66%  '
Set-TextValue $ws.Range('D32') '2.21'
$ws.Range('E32').Value = '  -1.00%  '
Set-TextValue $ws.Range('D33') '7.22'
$ws.Range('E33').Value = '  -2.84%  '
Set-TextValue $ws.Range('D34') '28.79'
$ws.Range('E34').Value = '  -2.21%  '
Set-TextValue $ws.Range('D35') '9.14'
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('E36').Value = '  -5.25%  '
$ws.Range('E37').Value = '  +0.23%  '
Set-TextValue $ws.Range('D38') '3.731.78'
$ws.Range('E38').Value = '  -1.61%  '
Set-TextValue $ws.Range('D39') '6.49'
$ws.Range('E39').Value = '  +8.27%  '
Set-TextValue $ws.Range('D40') '0.0985'
$ws.Range('E40').Value = '  -4.06%  '
Set-TextValue $ws.Range('D41') '2.32'
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('E42').Value = '  -3.13%  '
$ws.Range('E44').Value = '  +0.13%  '
Set-TextValue $ws.Range('D45') '3.09'
$ws.Range('E45').Value = '  -6.87%  '
Set-TextValue $ws.Range('D46') '162.16'
$ws.Range('E46').Value = '  -1.04%  '
Set-TextValue $ws.Range('D47') '48.79'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('E48').Value = '  -5.24%  '
$ws.Range('B49').Value = 'ONDO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range('D49') '1.37'
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D50') '408.38'
$ws.Range('E50').Value = '  -4.57%  '
Set-TextValue $ws.Range('D51') '8.54'
$ws.Range('E51').Value = '  -2.51%  '
